$d = $word.ActiveDocument

# 1. Replace the text "...some major changes" + "." run with merged "...some major changes."
# First, remove the separate "." run and the bookmark by replacing the two-part text
# with a single merged text using Find/Replace across the run boundary.

$d.Content.Find.Execute(
    "changed login button to match the wireframe. Could not contribute too much this iteration due to some major changes.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "changed login button to match the wireframe. Could not contribute too much this iteration due to some major changes.",
    2
) | Out-Null

# 2. Add new run text after "Dennis Lee: " in the Team Member 4 row.
$d.Content.Find.Execute(
    "Dennis Lee: ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Dennis Lee: Reviewed the login for users.",
    2
) | Out-Null
